$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Fill in "Day 7" (column H) answers for the most recent week
#    block (rows 44-57), which previously had empty H cells.
# ------------------------------------------------------------------
$ws.Range("H44").Value = "7：04"
$ws.Range("H45").Value = "7：24"
$ws.Range("H46").Value = "24：40"
$ws.Range("H47").Value = "24：40"
$ws.Range("H48").Value = 5
$ws.Range("H49").Value = 1
$ws.Range("H50").Value = 5
$ws.Range("H51").Value = 360
$ws.Range("H52").Value = "无"
$ws.Range("H53").Value = 20
$ws.Range("H54").Value = 3
$ws.Range("H55").Value = 2
$ws.Range("H56").Value = 2
$ws.Range("H57").Value = "无"

# ------------------------------------------------------------------
# 2) Append a brand-new (empty) week block in rows 59-76, built from
#    a full value+format copy of the previous week block (rows
#    40-57), which carries over the exact same cell styles.
# ------------------------------------------------------------------
$ws.Range("A40:H57").Copy()
$ws.Range("A59:H76").PasteSpecial(-4122)
$ws.Range("A40:H57").Copy()
$ws.Range("A59:H76").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Row heights for the newly-created rows (matches source rows 40-57).
$ws.Rows.Item(59).RowHeight = 22.5
$ws.Rows.Item(60).RowHeight = 17.25
$ws.Rows.Item(61).RowHeight = 16.5
$ws.Rows.Item(62).RowHeight = 54
$ws.Rows.Item(63).RowHeight = 33
$ws.Rows.Item(64).RowHeight = 16.5
$ws.Rows.Item(65).RowHeight = 16.5
$ws.Rows.Item(66).RowHeight = 16.5
$ws.Rows.Item(67).RowHeight = 49.5
$ws.Rows.Item(68).RowHeight = 16.5
$ws.Rows.Item(69).RowHeight = 33
$ws.Rows.Item(70).RowHeight = 33
$ws.Rows.Item(71).RowHeight = 82.5
$ws.Rows.Item(72).RowHeight = 66
$ws.Rows.Item(73).RowHeight = 49.5
$ws.Rows.Item(74).RowHeight = 99
$ws.Rows.Item(75).RowHeight = 115.5
$ws.Rows.Item(76).RowHeight = 66

# New week's 7 day/date headers (row 62).
$ws.Range("B62").Value = "第一天`n日期:2025-12-13"
$ws.Range("C62").Value = "第二天`n日期:2025-12-14"
$ws.Range("D62").Value = "第三天`n日期:2025-12-15"
$ws.Range("E62").Value = "第四天`n日期:2025-12-16"
$ws.Range("F62").Value = "第五天`n日期:2025-12-17"
$ws.Range("G62").Value = "第六天`n日期:2025-12-18"
$ws.Range("H62").Value = "第七天`n日期:2025-12-19"

# The new week's answer grid (rows 63-76, columns B:H) starts blank -
# clear the values that were copied along with the previous week's
# answers, while keeping the copied cell formatting intact.
$ws.Range("B63:H76").ClearContents()

# Merge the title/instruction/name rows of the new block, same as the
# other week blocks.
$ws.Range("A59:H59").Merge()
$ws.Range("A60:H60").Merge()
$ws.Range("B61:H61").Merge()

# ------------------------------------------------------------------
# 3) Update the sheet's active selection to match the new edit
#    position.
# ------------------------------------------------------------------
$ws.Range("M59").Select()
